# Add dummy data to excel files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student ID column (A) should be stored as text (numFmtId 49 -> "@" text format),
# applied down through the trailing blank rows too.
$ws.Range("A1:A16").NumberFormat = "@"

# Data entry order matches the shared-string insertion order seen in the
# target file: all of column A (Student ID) first, then row-by-row across
# columns B..E for rows 2-5.
$ws.Range("A2").Value = "0000001"
$ws.Range("A3").Value = "0000002"
$ws.Range("A4").Value = "0000003"
$ws.Range("A5").Value = "0000004"

# Row 2 - Doe, John, Jack, B
$ws.Range("B2").Value = "Doe"
$ws.Range("C2").Value = "John"
$ws.Range("D2").Value = "Jack"
$ws.Range("E2").Value = "B"

# Row 3 - Smith, James, (no preferred name), A
$ws.Range("B3").Value = "Smith"
$ws.Range("C3").Value = "James"
$ws.Range("E3").Value = "A"

# Row 4 - Patterson, Scott, (no preferred name), B
$ws.Range("B4").Value = "Patterson"
$ws.Range("C4").Value = "Scott "
$ws.Range("E4").Value = "B"

# Row 5 - Ditto, John, Jack, A
$ws.Range("B5").Value = "Ditto"
$ws.Range("C5").Value = "John"
$ws.Range("D5").Value = "Jack"
$ws.Range("E5").Value = "A"

# Column F - DOB. Copy the format from the existing date cell (F2) down first
# so the new cells reuse the existing date style instead of minting a new one.
$ws.Range("F2").Copy()
$ws.Range("F3:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F2").Value = "09/09/1998"
$ws.Range("F3").Value = "08/13/1999"
$ws.Range("F4").Value = "04/14/1998"
$ws.Range("F5").Value = "03/11/1997"

# Selection moves to D1 as in the final file
$ws.Range("D1").Select() | Out-Null
